$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 43

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 44656
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112030
$ws.Cells.Item($newRow, 7).Value = "Poroto granado"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 70
$ws.Cells.Item($newRow, 11).Value = 27000
$ws.Cells.Item($newRow, 12).Value = 27000
$ws.Cells.Item($newRow, 13).Value = 27000
$ws.Cells.Item($newRow, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 1080
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
